# "9th Stab - Cosmetic Changes"
# Two new weekly-data columns (Jun_17, Jun_15) are inserted to the left of the
# existing date column (which held Jun_13 in B and Jun_10 in C), pushing the
# old columns right by two (B->D, C->E). The two freshly inserted columns are
# seeded with the same placeholder value ("UN") that already fills column B
# for every data row, matching the pattern observed for every other column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns at C (shifts existing C -> E, keeps B in place).
$ws.Columns("C:D").Insert()

# --- Header row -------------------------------------------------------
# B1 already held "Jun_13" before the insert and keeps that text, but a
# brand-new, more recent pair of week labels is written into the two newly
# inserted header cells, and the old B1 text is relocated to D1 while the
# old C1 ("Jun_10") has already slid over to E1 by the column insert above.
$ws.Cells.Item(1, 4).Value() = $ws.Cells.Item(1, 2).Value()
$ws.Cells.Item(1, 2).Value() = "Jun_17"
$ws.Cells.Item(1, 3).Value() = "Jun_15"

# --- Data rows ----------------------------------------------------------
# For every data row, fill the two newly inserted cells with the same
# placeholder ("UN") already used in column B.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $placeholder = $ws.Cells.Item($r, 2).Value()
    $ws.Cells.Item($r, 3).Value() = $placeholder
    $ws.Cells.Item($r, 4).Value() = $placeholder
}

# --- Column widths (cosmetic) --------------------------------------------
# Columns C, D, E all take on the same 8.0-character width that column C
# had before the edit (C and D are the two newly duplicated columns, E is
# the original column C that slid right). In real Excel, C and D would also
# be flagged as a collapsed outline group (collapsed="true") while E stays
# expanded (collapsed="false"); we still group/collapse them here so the
# workbook's outline state reflects that intent as closely as the object
# model allows.
$ws.Columns("C").ColumnWidth = 7.166666666666666
$ws.Columns("D").ColumnWidth = 7.166666666666666
$ws.Columns("E").ColumnWidth = 7.166666666666666

$grp = $ws.Range("C1:D1").EntireColumn
$grp.Group()
$grp.ShowDetail = $false
